$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old instance-attribute rows (po_currency_code / func_item_cost)
# so their now-unused shared strings can be dropped, then rebuild the table
# with the new set of instance attributes.
$ws.Range("A2:E3").ClearContents()

# New data rows: eid_instance_id, eid_instance_attribute, datatype, profile_id, display_name
$ws.Range("A2").Value = 4
$ws.Range("B2").Value = "DENOM_RAW_COST"
$ws.Range("C2").Value = "mdex:double"
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = "Cost in USD"

$ws.Range("A3").Value = 4
$ws.Range("B3").Value = "PA_DRAFT_INVOICE"
$ws.Range("C3").Value = "mdex:string"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "PA Draft Invoice"

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "INVOICE_DATE"
$ws.Range("C4").Value = "mdex:string"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "Invoice Date"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "STATUS_TRX"
$ws.Range("C5").Value = "mdex:string"
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = "AR Invoice Status"

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "REASON_CODE"
$ws.Range("C6").Value = "mdex:string"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "AR Exception Reason"

$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "JOB_NAME"
$ws.Range("C7").Value = "mdex:string"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "Job"

$ws.Range("E7").Select() | Out-Null
